$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "fecha"
$ws.Cells.Item(2, 2).Value = "Sin categoría"
$ws.Cells.Item(3, 1).Value = "Saldo_de_la_deuda_externa_total_porcentaje_del_PIB_Deuda externa (pública y privada)_bond"
$ws.Cells.Item(3, 2).Value = "bond"
$ws.Cells.Item(4, 1).Value = "Reservas_internacionales_brutasDato_fin_de_mes_Brute_Internacional_Reservoir_economics"
$ws.Cells.Item(4, 2).Value = "economics"
$ws.Cells.Item(5, 1).Value = "Tasa_de_política_monetariaDato_diario_TPS_economics"
$ws.Cells.Item(5, 2).Value = "economics"
$ws.Cells.Item(6, 1).Value = "Índice_de_la_tasa_de_cambio_real_ITCR_según_IPP__B_ITCR_PP_exchange_rate"
$ws.Cells.Item(6, 2).Value = "exchange_rate"
$ws.Cells.Item(7, 1).Value = "Índice_de_precios_de_exportaciones_según_ippDato_f_Indice_Intercambios_IPP_index_pricing"
$ws.Cells.Item(7, 2).Value = "exports"
$ws.Cells.Item(8, 1).Value = "Índice_COLCAPDato_diario_Índice_COLCAP_index_pricing"
$ws.Cells.Item(8, 2).Value = "index_pricing"
$ws.Cells.Item(9, 1).Value = "Total_Exportaciones_Tradicionales_anex-EXPORTACIONES-SerieCafeCarbonPetroleoNotradicionales-mar2025_exports"
$ws.Cells.Item(9, 2).Value = "exports"
$ws.Cells.Item(10, 1).Value = "PRICE_Australia_10Y_Bond_bond"
$ws.Cells.Item(10, 2).Value = "bond"
$ws.Cells.Item(11, 1).Value = "PRICE_Italy_10Y_Bond_bond"
$ws.Cells.Item(11, 2).Value = "bond"
$ws.Cells.Item(12, 1).Value = "PRICE_Japan_10Y_Bond_bond"
$ws.Cells.Item(12, 2).Value = "bond"
$ws.Cells.Item(13, 1).Value = "PRICE_UK_10Y_Bond_bond"
$ws.Cells.Item(13, 2).Value = "bond"
$ws.Cells.Item(14, 1).Value = "PRICE_Germany_10Y_Bond_bond"
$ws.Cells.Item(14, 2).Value = "bond"
$ws.Cells.Item(15, 1).Value = "PRICE_Canada_10Y_Bond_bond"
$ws.Cells.Item(15, 2).Value = "bond"
$ws.Cells.Item(16, 1).Value = "PRICE_China_10Y_Bond_bond"
$ws.Cells.Item(16, 2).Value = "bond"
$ws.Cells.Item(17, 1).Value = "PRICE_CrudeOil_WTI_commodities"
$ws.Cells.Item(17, 2).Value = "commodities"
$ws.Cells.Item(18, 1).Value = "PRICE_Gold_Spot_commodities"
$ws.Cells.Item(18, 2).Value = "commodities"
$ws.Cells.Item(19, 1).Value = "PRICE_Silver_Spot_commodities"
$ws.Cells.Item(19, 2).Value = "commodities"
$ws.Cells.Item(20, 1).Value = "PRICE_Copper_Futures_commodities"
$ws.Cells.Item(20, 2).Value = "commodities"
$ws.Cells.Item(21, 1).Value = "PRICE_Platinum_Spot_commodities"
$ws.Cells.Item(21, 2).Value = "commodities"
$ws.Cells.Item(22, 1).Value = "PRICE_EUR_USD_Spot_exchange_rate"
$ws.Cells.Item(22, 2).Value = "exchange_rate"
$ws.Cells.Item(23, 1).Value = "PRICE_GBP_USD_Spot_exchange_rate"
$ws.Cells.Item(23, 2).Value = "exchange_rate"
$ws.Cells.Item(24, 1).Value = "PRICE_JPY_USD_Spot_exchange_rate"
$ws.Cells.Item(24, 2).Value = "exchange_rate"
$ws.Cells.Item(25, 1).Value = "PRICE_CNY_USD_Spot_exchange_rate"
$ws.Cells.Item(25, 2).Value = "exchange_rate"
$ws.Cells.Item(26, 1).Value = "PRICE_AUD_USD_Spot_exchange_rate"
$ws.Cells.Item(26, 2).Value = "exchange_rate"
$ws.Cells.Item(27, 1).Value = "PRICE_CAD_USD_Spot_exchange_rate"
$ws.Cells.Item(27, 2).Value = "exchange_rate"
$ws.Cells.Item(28, 1).Value = "PRICE_MXN_USD_Spot_exchange_rate"
$ws.Cells.Item(28, 2).Value = "exchange_rate"
$ws.Cells.Item(29, 1).Value = "PRICE_EUR_GBP_Cross_exchange_rate"
$ws.Cells.Item(29, 2).Value = "exchange_rate"
$ws.Cells.Item(30, 1).Value = "PRICE_S&P500_Index_index_pricing"
$ws.Cells.Item(30, 2).Value = "index_pricing"
$ws.Cells.Item(31, 1).Value = "PRICE_NASDAQ_Composite_index_pricing"
$ws.Cells.Item(31, 2).Value = "index_pricing"
$ws.Cells.Item(32, 1).Value = "PRICE_Russell_2000_index_pricing"
$ws.Cells.Item(32, 2).Value = "index_pricing"
$ws.Cells.Item(33, 1).Value = "PRICE_FTSE_100_index_pricing"
$ws.Cells.Item(33, 2).Value = "index_pricing"
$ws.Cells.Item(34, 1).Value = "PRICE_Nikkei_225_index_pricing"
$ws.Cells.Item(34, 2).Value = "index_pricing"
$ws.Cells.Item(35, 1).Value = "PRICE_DAX_30_index_pricing"
$ws.Cells.Item(35, 2).Value = "index_pricing"
$ws.Cells.Item(36, 1).Value = "PRICE_Shanghai_Composite_index_pricing"
$ws.Cells.Item(36, 2).Value = "index_pricing"
$ws.Cells.Item(37, 1).Value = "PRICE_VIX_VolatilityIndex_index_pricing"
$ws.Cells.Item(37, 2).Value = "index_pricing"
$ws.Cells.Item(38, 1).Value = "PRICE_Aluminium_Spot_commodities"
$ws.Cells.Item(38, 2).Value = "commodities"
$ws.Cells.Item(39, 1).Value = "PRICE_Lead_Futures_commodities"
$ws.Cells.Item(39, 2).Value = "commodities"
$ws.Cells.Item(40, 1).Value = "PRICE_Palladium_Futures_commodities"
$ws.Cells.Item(40, 2).Value = "commodities"
$ws.Cells.Item(41, 1).Value = "PRICE_Tin_Futures_commodities"
$ws.Cells.Item(41, 2).Value = "commodities"
$ws.Cells.Item(42, 1).Value = "PRICE_Zinc_Futures_commodities"
$ws.Cells.Item(42, 2).Value = "commodities"
$ws.Cells.Item(43, 1).Value = "PRICE_Nickel_Futures_commodities"
$ws.Cells.Item(43, 2).Value = "commodities"
$ws.Cells.Item(44, 1).Value = "PRICE_Iron_ore_fines_62%_Fe_CFR_Futures_commodities"
$ws.Cells.Item(44, 2).Value = "commodities"
$ws.Cells.Item(45, 1).Value = "PRICE_Brent_Oil_Futures_commodities"
$ws.Cells.Item(45, 2).Value = "commodities"
$ws.Cells.Item(46, 1).Value = "PRICE_Dow_Jones_Industrial_Average_index_pricing"
$ws.Cells.Item(46, 2).Value = "index_pricing"
$ws.Cells.Item(47, 1).Value = "PRICE_Nasdaq_100_Fi_index_pricing"
$ws.Cells.Item(47, 2).Value = "index_pricing"
$ws.Cells.Item(48, 1).Value = "PRICE_CAC_40_index_pricing"
$ws.Cells.Item(48, 2).Value = "index_pricing"
$ws.Cells.Item(49, 1).Value = "PRICE_JPY_EUR_exchange_rate"
$ws.Cells.Item(49, 2).Value = "exchange_rate"
$ws.Cells.Item(50, 1).Value = "PRICE_BRL_NZD_exchange_rate"
$ws.Cells.Item(50, 2).Value = "exchange_rate"
$ws.Cells.Item(51, 1).Value = "PRICE_NOK_AUD_exchange_rate"
$ws.Cells.Item(51, 2).Value = "exchange_rate"
$ws.Cells.Item(52, 1).Value = "PRICE_PHP_ZAR_exchange_rate"
$ws.Cells.Item(52, 2).Value = "exchange_rate"
$ws.Cells.Item(53, 1).Value = "PRICE_USDCOP-US_Dollar_Colombian_Peso_exchange_rate"
$ws.Cells.Item(53, 2).Value = "exchange_rate"
$ws.Cells.Item(54, 1).Value = "PRICE_Colombia_5_Year_Bond_bond"
$ws.Cells.Item(54, 2).Value = "bond"
$ws.Cells.Item(55, 1).Value = "PRICE_Colombia_10_Year_Bond_bond"
$ws.Cells.Item(55, 2).Value = "bond"
$ws.Cells.Item(56, 1).Value = "ESI_GACDISA_US_Empire_State_Index_business_confidence"
$ws.Cells.Item(56, 2).Value = "business_confidence"
$ws.Cells.Item(57, 1).Value = "ESI_AWCDISA_US_Empire_State_Index_business_confidence"
$ws.Cells.Item(57, 2).Value = "business_confidence"
$ws.Cells.Item(58, 1).Value = "Put_strike_Put_Call_Ratio_SPY_consumer_confidence"
$ws.Cells.Item(58, 2).Value = "consumer_confidence"
$ws.Cells.Item(59, 1).Value = "Put_bid_Put_Call_Ratio_SPY_consumer_confidence"
$ws.Cells.Item(59, 2).Value = "consumer_confidence"
$ws.Cells.Item(60, 1).Value = "Put_ask_Put_Call_Ratio_SPY_consumer_confidence"
$ws.Cells.Item(60, 2).Value = "consumer_confidence"
$ws.Cells.Item(61, 1).Value = "Put_vol_Put_Call_Ratio_SPY_consumer_confidence"
$ws.Cells.Item(61, 2).Value = "consumer_confidence"
$ws.Cells.Item(62, 1).Value = "Put_delta_Put_Call_Ratio_SPY_consumer_confidence"
$ws.Cells.Item(62, 2).Value = "consumer_confidence"
$ws.Cells.Item(63, 1).Value = "Put_gamma_Put_Call_Ratio_SPY_consumer_confidence"
$ws.Cells.Item(63, 2).Value = "consumer_confidence"
$ws.Cells.Item(64, 1).Value = "Put_theta_Put_Call_Ratio_SPY_consumer_confidence"
$ws.Cells.Item(64, 2).Value = "consumer_confidence"
$ws.Cells.Item(65, 1).Value = "Put_vega_Put_Call_Ratio_SPY_consumer_confidence"
$ws.Cells.Item(65, 2).Value = "consumer_confidence"
$ws.Cells.Item(66, 1).Value = "Put_rho_Put_Call_Ratio_SPY_consumer_confidence"
$ws.Cells.Item(66, 2).Value = "consumer_confidence"
$ws.Cells.Item(67, 1).Value = "NFCI_Chicago_Fed_NFCI_leading_economic_index"
$ws.Cells.Item(67, 2).Value = "economics"
$ws.Cells.Item(68, 1).Value = "ANFCI_Chicago_Fed_NFCI_leading_economic_index"
$ws.Cells.Item(68, 2).Value = "economics"
$ws.Cells.Item(69, 1).Value = "TasasTES_Tasa de interés Cero Cupón, Títulos de Tesorería (TES), pesos - 1 año_TasasTES_Bond_bond"
$ws.Cells.Item(69, 2).Value = "bond"
$ws.Cells.Item(70, 1).Value = "TasasTES_Tasa de interés Cero Cupón, Títulos de Tesorería (TES), pesos - 5 años_TasasTES_Bond_bond"
$ws.Cells.Item(70, 2).Value = "bond"
$ws.Cells.Item(71, 1).Value = "TasasTES_Tasa de interés Cero Cupón, Títulos de Tesorería (TES), pesos - 10 años_TasasTES_Bond_bond"
$ws.Cells.Item(71, 2).Value = "bond"
$ws.Cells.Item(72, 1).Value = "TasasTES_Tasa de interés Cero Cupón, Títulos de Tesorería (TES), UVR - 1 año_TasasTES_Bond_bond"
$ws.Cells.Item(72, 2).Value = "bond"
$ws.Cells.Item(73, 1).Value = "TasasTES_Tasa de interés Cero Cupón, Títulos de Tesorería (TES), UVR - 5 años_TasasTES_Bond_bond"
$ws.Cells.Item(73, 2).Value = "bond"
$ws.Cells.Item(74, 1).Value = "TasasTES_Tasa de interés Cero Cupón, Títulos de Tesorería (TES), UVR - 10 años_TasasTES_Bond_bond"
$ws.Cells.Item(74, 2).Value = "bond"
$ws.Cells.Item(75, 1).Value = "Actual_US_ISM_Manufacturing_business_confidence"
$ws.Cells.Item(75, 2).Value = "business_confidence"
$ws.Cells.Item(76, 1).Value = "Actual_US_ISM_Services_business_confidence"
$ws.Cells.Item(76, 2).Value = "business_confidence"
$ws.Cells.Item(77, 1).Value = "Actual_US_Philly_Fed_Index_business_confidence"
$ws.Cells.Item(77, 2).Value = "business_confidence"
$ws.Cells.Item(78, 1).Value = "Actual_France_Business_Climate_business_confidence"
$ws.Cells.Item(78, 2).Value = "business_confidence"
$ws.Cells.Item(79, 1).Value = "Actual_EuroZone_Business_Climate_business_confidence"
$ws.Cells.Item(79, 2).Value = "business_confidence"
$ws.Cells.Item(80, 1).Value = "Actual_U.S. All Car Sales_car_registrations"
$ws.Cells.Item(80, 2).Value = "car_registrations"
$ws.Cells.Item(81, 1).Value = "Actual_US_Consumer_Confidence_consumer_confidence"
$ws.Cells.Item(81, 2).Value = "consumer_confidence"
$ws.Cells.Item(82, 1).Value = "Actual_China_PMI_Manufacturing_economics"
$ws.Cells.Item(82, 2).Value = "economics"
$ws.Cells.Item(83, 1).Value = "Actual_Singapore_NonOil_Exports_YoY_economics"
$ws.Cells.Item(83, 2).Value = "exports"
$ws.Cells.Item(84, 1).Value = "Actual_Japan_M2_MoneySupply_YoY_economics"
$ws.Cells.Item(84, 2).Value = "economics"
$ws.Cells.Item(85, 1).Value = "Actual_China_M2_MoneySupply_YoY_economics"
$ws.Cells.Item(85, 2).Value = "economics"
$ws.Cells.Item(86, 1).Value = "Actual_Mexico_CPI_MoM_economics"
$ws.Cells.Item(86, 2).Value = "economics"
$ws.Cells.Item(87, 1).Value = "Actual_BOJ_Policy_Rate_economics"
$ws.Cells.Item(87, 2).Value = "economics"
$ws.Cells.Item(88, 1).Value = "Actual_UK_Retail_Sales_MoM_economics"
$ws.Cells.Item(88, 2).Value = "economics"
$ws.Cells.Item(89, 1).Value = "Actual_China_Exports_exports"
$ws.Cells.Item(89, 2).Value = "exports"
$ws.Cells.Item(90, 1).Value = "Actual_US_Exports_exports"
$ws.Cells.Item(90, 2).Value = "exports"
$ws.Cells.Item(91, 1).Value = "Actual_US_ConferenceBoard_LEI_leading_economic_index"
$ws.Cells.Item(91, 2).Value = "economics"
$ws.Cells.Item(92, 1).Value = "Actual_Japan_Leading_Indicator_leading_economic_index"
$ws.Cells.Item(92, 2).Value = "economics"
$ws.Cells.Item(93, 1).Value = "Actual_China_Unemployment_Rate_unemployment_rate"
$ws.Cells.Item(93, 2).Value = "unemployment_rate"
$ws.Cells.Item(94, 1).Value = "Actual_Eurozone_Unemployment_Rate_unemployment_rate"
$ws.Cells.Item(94, 2).Value = "unemployment_rate"
$ws.Cells.Item(95, 1).Value = "Actual_IPC Colombia_economics"
$ws.Cells.Item(95, 2).Value = "economics"
$ws.Cells.Item(96, 1).Value = "Actual_IPC EEUU_economics"
$ws.Cells.Item(96, 2).Value = "economics"
$ws.Cells.Item(97, 1).Value = "DGS10_US_10Y_Treasury_bond"
$ws.Cells.Item(97, 2).Value = "bond"
$ws.Cells.Item(98, 1).Value = "DGS2_US_2Y_Treasury_bond"
$ws.Cells.Item(98, 2).Value = "bond"
$ws.Cells.Item(99, 1).Value = "AAA_Corporate_Bond_AAA_Spread_bond"
$ws.Cells.Item(99, 2).Value = "bond"
$ws.Cells.Item(100, 1).Value = "BAA10YM_Corporate_Bond_BBB_Spread_bond"
$ws.Cells.Item(100, 2).Value = "bond"
$ws.Cells.Item(101, 1).Value = "BAMLH0A0HYM2_High_Yield_Bond_Spread_bond"
$ws.Cells.Item(101, 2).Value = "economics"
$ws.Cells.Item(102, 1).Value = "DNKSLRTCR03GPSAM_Denmark_Car_Registrations_MoM_car_registrations"
$ws.Cells.Item(102, 2).Value = "car_registrations"
$ws.Cells.Item(103, 1).Value = "USASLRTCR03GPSAM_US_Car_Registrations_MoM_car_registrations"
$ws.Cells.Item(103, 2).Value = "car_registrations"
$ws.Cells.Item(104, 1).Value = "ZAFSLRTCR03GPSAM_SouthAfrica_Car_Registrations_MoM_car_registrations"
$ws.Cells.Item(104, 2).Value = "car_registrations"
$ws.Cells.Item(105, 1).Value = "GBRSLRTCR03GPSAM_United_Kingdom_Car_Registrations_MoM_car_registrations"
$ws.Cells.Item(105, 2).Value = "car_registrations"
$ws.Cells.Item(106, 1).Value = "ESPSLRTCR03GPSAM_Spain_Car_Registrations_MoM_car_registrations"
$ws.Cells.Item(106, 2).Value = "car_registrations"
$ws.Cells.Item(107, 1).Value = "BUSLOANS_US_Commercial_Loans_comm_loans"
$ws.Cells.Item(107, 2).Value = "comm_loans"
$ws.Cells.Item(108, 1).Value = "CREACBM027NBOG_US_RealEstate_Commercial_Loans_comm_loans"
$ws.Cells.Item(108, 2).Value = "economics"
$ws.Cells.Item(109, 1).Value = "TOTALSL_US_Consumer_Credit_comm_loans"
$ws.Cells.Item(109, 2).Value = "comm_loans"
$ws.Cells.Item(110, 1).Value = "CSCICP02EZM460S_EuroZone_Consumer_Confidence_consumer_confidence"
$ws.Cells.Item(110, 2).Value = "economics"
$ws.Cells.Item(111, 1).Value = "CSCICP02CHQ460S_Switzerland_Consumer_Confidence_consumer_confidence"
$ws.Cells.Item(111, 2).Value = "consumer_confidence"
$ws.Cells.Item(112, 1).Value = "UMCSENT_Michigan_Consumer_Sentiment_consumer_confidence"
$ws.Cells.Item(112, 2).Value = "consumer_confidence"
$ws.Cells.Item(113, 1).Value = "CPIAUCSL_US_CPI_economics"
$ws.Cells.Item(113, 2).Value = "economics"
$ws.Cells.Item(114, 1).Value = "CPILFESL_US_Core_CPI_economics"
$ws.Cells.Item(114, 2).Value = "economics"
$ws.Cells.Item(115, 1).Value = "PCE_US_PCE_economics"
$ws.Cells.Item(115, 2).Value = "economics"
$ws.Cells.Item(116, 1).Value = "PCEPILFE_US_Core_PCE_economics"
$ws.Cells.Item(116, 2).Value = "economics"
$ws.Cells.Item(117, 1).Value = "PPIACO_US_PPI_economics"
$ws.Cells.Item(117, 2).Value = "economics"
$ws.Cells.Item(118, 1).Value = "INDPRO_US_Industrial_Production_MoM_economics"
$ws.Cells.Item(118, 2).Value = "economics"
$ws.Cells.Item(119, 1).Value = "CSUSHPINSA_US_CaseShiller_HomePrice_economics"
$ws.Cells.Item(119, 2).Value = "economics"
$ws.Cells.Item(120, 1).Value = "GDP_US_GDP_Growth_economics"
$ws.Cells.Item(120, 2).Value = "economics"
$ws.Cells.Item(121, 1).Value = "TCU_US_Capacity_Utilization_economics"
$ws.Cells.Item(121, 2).Value = "economics"
$ws.Cells.Item(122, 1).Value = "PERMIT_US_Building_Permits_economics"
$ws.Cells.Item(122, 2).Value = "economics"
$ws.Cells.Item(123, 1).Value = "HOUST_US_Housing_Starts_economics"
$ws.Cells.Item(123, 2).Value = "economics"
$ws.Cells.Item(124, 1).Value = "FEDFUNDS_US_FedFunds_Rate_economics"
$ws.Cells.Item(124, 2).Value = "economics"
$ws.Cells.Item(125, 1).Value = "ECBDFR_ECB_Deposit_Rate_economics"
$ws.Cells.Item(125, 2).Value = "economics"
$ws.Cells.Item(126, 1).Value = "WALCL_Fed_Balance_Sheet_economics"
$ws.Cells.Item(126, 2).Value = "economics"
$ws.Cells.Item(127, 1).Value = "DTWEXBGS_Dollar_Index_DXY_index_pricing"
$ws.Cells.Item(127, 2).Value = "index_pricing"
$ws.Cells.Item(128, 1).Value = "UNRATE_US_Unemployment_Rate_unemployment_rate"
$ws.Cells.Item(128, 2).Value = "unemployment_rate"
$ws.Cells.Item(129, 1).Value = "PAYEMS_US_Nonfarm_Payrolls_unemployment_rate"
$ws.Cells.Item(129, 2).Value = "unemployment_rate"
$ws.Cells.Item(130, 1).Value = "ICSA_US_Initial_Jobless_Claims_unemployment_rate"
$ws.Cells.Item(130, 2).Value = "unemployment_rate"
$ws.Cells.Item(131, 1).Value = "DGS10_DGS10_bond"
$ws.Cells.Item(131, 2).Value = "bond"
